# Update the "Groups" sheet to use the new bracket-range well-selector syntax
# instead of explicit comma-separated well lists.
$wb = $excel.ActiveWorkbook
$groups = $wb.Worksheets.Item("Groups")

$groups.Range("B4").Value = "plate_01_A[1:3]"
$groups.Range("B3").Value = "plate_01_A[3,8:-1:7]"
$groups.Range("B2").Value = "plate_01_A[1,5,9]"

# Reflect the new active sheet/selection recorded in the workbook view state:
# the "Groups" sheet becomes the active tab (previously "Views"), with F11
# selected there; "Views" keeps its own last selection (G22) but is no
# longer the tabSelected sheet.
$groups.Activate()
$groups.Range("F11").Select()
